# Applies the "Updated symbol list" GitHub Actions refresh (Fri Feb 10 2023):
# rows 6-17 rotate by one slot (a new GateToken row is inserted at the top of
# that block, pushing FTXToken..LEO down one row and dropping the prior last
# slot), every row's Volume(1h)/Hora are refreshed, and many Price cells move
# by small live-market deltas.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="305.79"; E="-0.53%"; G="23" },
    @{ Row=3; D="40.44"; E="3.34%"; G="23" },
    @{ Row=4; D="5.110"; E="1.55%"; G="23" },
    @{ Row=5; D="0.07585"; E="-2.05%"; G="23" },
    @{ Row=6; B="GateToken"; C="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D="4.270"; E="-0.60%"; G="23" },
    @{ Row=7; B="FTXToken"; C="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D="1.631"; E="0.79%"; G="23" },
    @{ Row=8; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="0.9058"; E="-1.52%"; G="23" },
    @{ Row=9; B="BTSEToken"; C="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D="2.423"; E="-6.08%"; G="23" },
    @{ Row=10; B="LiechtensteinCryptoassetsExchange"; C="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D="0.1010"; E="2.23%"; G="23" },
    @{ Row=11; B="WazirX"; C="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D="0.1746"; E="1.31%"; G="23" },
    @{ Row=12; B="MandalaExchangeToken"; C="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D="0.09112"; E="2.60%"; G="23" },
    @{ Row=13; B="BitrueCoin"; C="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D="0.04179"; E="-5.14%"; G="23" },
    @{ Row=14; B="BitMartToken"; C="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D="0.1056"; E="-0.40%"; G="23" },
    @{ Row=15; B="BitForexToken"; C="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D="0.001255"; E="0.44%"; G="23" },
    @{ Row=16; B="TigerCash"; C="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D="0.005858"; E="3.35%"; G="23" },
    @{ Row=17; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.352"; E="-0.40%"; G="23" },
    @{ Row=18; E="-2.81%"; G="23" },
    @{ Row=19; D="6.632"; E="-5.76%"; G="23" },
    @{ Row=20; E="-0.84%"; G="23" },
    @{ Row=21; D="0.2729"; E="-2.00%"; G="23" },
    @{ Row=22; D="0.04177"; E="1.20%"; G="23" },
    @{ Row=23; D="0.001229"; E="1.91%"; G="23" },
    @{ Row=24; D="0.004056"; E="-0.60%"; G="23" },
    @{ Row=25; D="0.0001302"; E="6.31%"; G="23" },
    @{ Row=26; D="0.0003012"; E="0.35%"; G="23" },
    @{ Row=27; G="23" },
    @{ Row=28; G="23" },
    @{ Row=29; G="23" },
    @{ Row=30; G="23" },
    @{ Row=31; G="23" },
    @{ Row=32; G="23" },
    @{ Row=33; G="23" },
    @{ Row=34; G="23" },
    @{ Row=35; G="23" },
    @{ Row=36; G="23" },
    @{ Row=37; G="23" },
    @{ Row=38; D="0.02367"; E="0.84%"; G="23" },
    @{ Row=39; D="0.05134"; E="-0.01%"; G="23" },
    @{ Row=40; D="0.007779"; E="-2.83%"; G="23" },
    @{ Row=41; D="0.1292"; E="-2.74%"; G="23" },
    @{ Row=42; D="0.007077"; E="-4.44%"; G="23" },
    @{ Row=43; E="-4.13%"; G="23" },
    @{ Row=44; D="0.008454"; E="16.12%"; G="23" },
    @{ Row=45; D="0.3332"; E="0.38%"; G="23" },
    @{ Row=46; D="0.00006354"; E="-4.43%"; G="23" },
    @{ Row=47; E="-0.59%"; G="23" },
    @{ Row=48; D="0.004408"; E="6.68%"; G="23" },
    @{ Row=49; D="0.007047"; E="108.22%"; G="23" },
    @{ Row=50; D="0.00002104"; E="-0.59%"; G="23" },
    @{ Row=51; D="0.0002004"; E="-0.59%"; G="23" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
    if ($u.ContainsKey("G")) {
        $ws.Cells.Item($u.Row, 7).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
}
